# Apply the "Add files via upload" edit:
#  - sheet1 (2025-03-18): append a "Total Tickets" row with a SUM formula
#  - sheet2 (2025-03-19): replace the day's data with new figures and append
#    a "Total Tickets" row (plain value)
#  - add a brand-new sheet "2025-03-20" with its own data and the same kind
#    of "Total Tickets" row, and make it the active sheet/tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "2025-03-18" - unchanged data, new totals row with formula
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("C8").Value = "Total Tickets"
$ws1.Range("D8").Formula = "=SUM(D2:D7)"
$ws1.Range("I10").Select()

# ---------------------------------------------------------------------
# Sheet 2: "2025-03-19" - new data values + totals row
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(2,1).Value = "Longside Upper"
$ws2.Cells.Item(2,2).Value = 83
$ws2.Cells.Item(2,3).Value = 218
$ws2.Cells.Item(2,4).Value = 86

$ws2.Cells.Item(3,1).Value = "Shortside Upper"
$ws2.Cells.Item(3,2).Value = 84
$ws2.Cells.Item(3,3).Value = 175
$ws2.Cells.Item(3,4).Value = 93

$ws2.Cells.Item(4,1).Value = "Longside Lower"
$ws2.Cells.Item(4,2).Value = 98
$ws2.Cells.Item(4,3).Value = 275
$ws2.Cells.Item(4,4).Value = 141

$ws2.Cells.Item(5,1).Value = "Shortside Lower"
$ws2.Cells.Item(5,2).Value = 100
$ws2.Cells.Item(5,3).Value = 218
$ws2.Cells.Item(5,4).Value = 99

$ws2.Cells.Item(6,1).Value = "Club Level"
$ws2.Cells.Item(6,2).Value = 249
$ws2.Cells.Item(6,3).Value = 689
$ws2.Cells.Item(6,4).Value = 24

$ws2.Cells.Item(7,1).Value = "VIP & Executive Box"
$ws2.Cells.Item(7,2).Value = 995
$ws2.Cells.Item(7,3).Value = 995
$ws2.Cells.Item(7,4).Value = 1

$ws2.Range("C8").Value = "Total Tickets"
$ws2.Range("D8").Value = 444

$ws2.Columns("A:D").ClearFormats()

$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# Sheet 3: new "2025-03-20" sheet appended after the last sheet
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "2025-03-20"

$ws3.Range("A1").Value = "Seat Type"
$ws3.Range("B1").Value = "Min_Price"
$ws3.Range("C1").Value = "Avg_Price"
$ws3.Range("D1").Value = "Ticket_Count"

$ws3.Cells.Item(2,1).Value = "Shortside Upper"
$ws3.Cells.Item(2,2).Value = 84
$ws3.Cells.Item(2,3).Value = 176
$ws3.Cells.Item(2,4).Value = 92

$ws3.Cells.Item(3,1).Value = "Longside Lower"
$ws3.Cells.Item(3,2).Value = 98
$ws3.Cells.Item(3,3).Value = 277
$ws3.Cells.Item(3,4).Value = 139

$ws3.Cells.Item(4,1).Value = "Shortside Lower"
$ws3.Cells.Item(4,2).Value = 100
$ws3.Cells.Item(4,3).Value = 218
$ws3.Cells.Item(4,4).Value = 98

$ws3.Cells.Item(5,1).Value = "Longside Upper"
$ws3.Cells.Item(5,2).Value = 100
$ws3.Cells.Item(5,3).Value = 220
$ws3.Cells.Item(5,4).Value = 85

$ws3.Cells.Item(6,1).Value = "Club Level"
$ws3.Cells.Item(6,2).Value = 249
$ws3.Cells.Item(6,3).Value = 668
$ws3.Cells.Item(6,4).Value = 20

$ws3.Cells.Item(7,1).Value = "VIP & Executive Box"
$ws3.Cells.Item(7,2).Value = 995
$ws3.Cells.Item(7,3).Value = 995
$ws3.Cells.Item(7,4).Value = 1

$ws3.Range("C8").Value = "Total Tickets"
$ws3.Range("D8").Value = 435

$ws3.Activate()
